# cryptos.xlsx refresh — GitHub Actions scheduled update (Sun Sep 15 23:42:10 UTC 2024)
#
# The upstream Coinranking-backed "cryptos list" sheet gets re-pulled on a
# cron job; each run rewrites the Price (D) and Volume(1h) (E) columns with
# freshly scraped values for every coin row, and occasionally a couple of
# rows trade places (their relative 1h change crossed over), which also
# swaps that row's Coin name (B) and Link (C).
#
# All of Price/Volume/Coin/Link are stored as *text* in this sheet (plain
# strings, not numbers/percentages), so values like "553.36" or "1.00" must
# stay literal text instead of being auto-converted to numbers by Excel's
# normal cell-entry parsing (which would e.g. turn "1.00" into 1 and drop
# the trailing zero). Cells whose new text looks number-like are written
# with a leading quote (forces text entry, same as typing '553.36 into
# Excel) and then have their style reset to "Normal" so no stray
# number-format/quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '59.114.85' },
    @{ Cell = 'E2'; Value = '  -1.44%  ' },
    @{ Cell = 'D3'; Value = '2.320.28' },
    @{ Cell = 'E3'; Value = '  -4.03%  ' },
    @{ Cell = 'E4'; Value = '  +0.10%  ' },
    @{ Cell = 'D5'; Value = '553.36' },
    @{ Cell = 'E5'; Value = '  +0.04%  ' },
    @{ Cell = 'D6'; Value = '131.51' },
    @{ Cell = 'E6'; Value = '  -4.11%  ' },
    @{ Cell = 'E7'; Value = '  +0.17%  ' },
    @{ Cell = 'D8'; Value = '0.571' },
    @{ Cell = 'E8'; Value = '  -4.42%  ' },
    @{ Cell = 'E9'; Value = '  -2.60%  ' },
    @{ Cell = 'D10'; Value = '5.54' },
    @{ Cell = 'E10'; Value = '  -2.32%  ' },
    @{ Cell = 'E11'; Value = '  +0.89%  ' },
    @{ Cell = 'D12'; Value = '0.338' },
    @{ Cell = 'E12'; Value = '  -4.63%  ' },
    @{ Cell = 'D13'; Value = '23.82' },
    @{ Cell = 'E13'; Value = '  -5.95%  ' },
    @{ Cell = 'D14'; Value = '2.739.72' },
    @{ Cell = 'E14'; Value = '  -3.83%  ' },
    @{ Cell = 'D15'; Value = '59.090.56' },
    @{ Cell = 'E15'; Value = '  -1.39%  ' },
    @{ Cell = 'E16'; Value = '  -2.53%  ' },
    @{ Cell = 'D17'; Value = '2.322.50' },
    @{ Cell = 'E17'; Value = '  -4.15%  ' },
    @{ Cell = 'D18'; Value = '10.82' },
    @{ Cell = 'E18'; Value = '  -4.13%  ' },
    @{ Cell = 'D19'; Value = '4.40' },
    @{ Cell = 'E19'; Value = '  -0.95%  ' },
    @{ Cell = 'D20'; Value = '316.96' },
    @{ Cell = 'E20'; Value = '  -3.31%  ' },
    @{ Cell = 'D21'; Value = '6.53' },
    @{ Cell = 'E21'; Value = '  -2.14%  ' },
    @{ Cell = 'D22'; Value = '0.999' },
    @{ Cell = 'E22'; Value = '  -0.05%  ' },
    @{ Cell = 'D23'; Value = '63.42' },
    @{ Cell = 'E23'; Value = '  -3.97%  ' },
    @{ Cell = 'E24'; Value = '  -4.04%  ' },
    @{ Cell = 'E25'; Value = '  +0.04%  ' },
    @{ Cell = 'D26'; Value = '8.33' },
    @{ Cell = 'E26'; Value = '  -3.55%  ' },
    @{ Cell = 'E27'; Value = '  -5.08%  ' },
    @{ Cell = 'E28'; Value = '  +0.29%  ' },
    @{ Cell = 'D29'; Value = '170.18' },
    @{ Cell = 'E29'; Value = '  +0.89%  ' },
    @{ Cell = 'D30'; Value = '0.0₃0738' },
    @{ Cell = 'E30'; Value = '  -4.94%  ' },
    @{ Cell = 'D31'; Value = '5.87' },
    @{ Cell = 'E31'; Value = '  -3.07%  ' },
    @{ Cell = 'D32'; Value = '1.08' },
    @{ Cell = 'E32'; Value = '  +4.21%  ' },
    @{ Cell = 'D33'; Value = '0.392' },
    @{ Cell = 'E33'; Value = '  -3.65%  ' },
    @{ Cell = 'E34'; Value = '  +0.03%  ' },
    @{ Cell = 'D35'; Value = '17.84' },
    @{ Cell = 'E35'; Value = '  -4.06%  ' },
    @{ Cell = 'B36'; Value = 'FirstDigitalUSD' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' },
    @{ Cell = 'D36'; Value = '1.00' },
    @{ Cell = 'E36'; Value = '  -0.03%  ' },
    @{ Cell = 'B37'; Value = 'ImmutableX' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Cell = 'D37'; Value = '1.28' },
    @{ Cell = 'E37'; Value = '  -2.79%  ' },
    @{ Cell = 'D38'; Value = '4.03' },
    @{ Cell = 'E38'; Value = '  -4.04%  ' },
    @{ Cell = 'D39'; Value = '1.55' },
    @{ Cell = 'E39'; Value = '  -3.37%  ' },
    @{ Cell = 'D40'; Value = '38.49' },
    @{ Cell = 'E40'; Value = '  -2.54%  ' },
    @{ Cell = 'D41'; Value = '304.21' },
    @{ Cell = 'E41'; Value = '  -7.58%  ' },
    @{ Cell = 'D42'; Value = '143.71' },
    @{ Cell = 'E42'; Value = '  +2.48%  ' },
    @{ Cell = 'D43'; Value = '3.46' },
    @{ Cell = 'E43'; Value = '  -5.55%  ' },
    @{ Cell = 'D44'; Value = '0.0952' },
    @{ Cell = 'E44'; Value = '  -2.10%  ' },
    @{ Cell = 'D45'; Value = '0.0500' },
    @{ Cell = 'E45'; Value = '  -3.24%  ' },
    @{ Cell = 'B46'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D46'; Value = '18.72' },
    @{ Cell = 'E46'; Value = '  -4.78%  ' },
    @{ Cell = 'B47'; Value = 'Mantle' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' },
    @{ Cell = 'D47'; Value = '0.558' },
    @{ Cell = 'E47'; Value = '  -3.46%  ' },
    @{ Cell = 'E48'; Value = '  -4.29%  ' },
    @{ Cell = 'D49'; Value = '11.04' },
    @{ Cell = 'E49'; Value = '  -0.01%  ' },
    @{ Cell = 'E50'; Value = '  +0.18%  ' },
    @{ Cell = 'E51'; Value = '  -0.66%  ' }
)

foreach ($update in $updates) {
    $cellRef = $update.Cell
    $newValue = $update.Value
    $range = $ws.Range($cellRef)

    # Does the new text look like a plain number (e.g. "553.36", "1.00",
    # "-0.05")? If so, Excel's default Value-assignment would silently
    # coerce it to a numeric cell and normalize its formatting. Force it to
    # stay text by entering it the way a user would to keep text that looks
    # like a number: a leading apostrophe.
    $trimmed = $newValue.Trim()
    $looksNumeric = $trimmed -match '^[-+]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $range.Value = "'" + $newValue
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}
